$wb = $excel.ActiveWorkbook

$wsAData = $wb.Worksheets.Item("Forecast Model A - Data")
$wsBData = $wb.Worksheets.Item("Forecast Model B - Data")

# Update the J-column (date) values on both data sheets.
# Rows 2-26 and rows 61-85 both get the same sequential run of dates
# (45991 .. 46015).
$rowBlocks = @(
    ,@(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26)
    ,@(61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85)
)

foreach ($sheet in @($wsAData, $wsBData)) {
    foreach ($rows in $rowBlocks) {
        $val = 45991
        foreach ($r in $rows) {
            $sheet.Cells.Item($r, 10).Value = $val
            $val = $val + 1
        }
    }
}

# Update the selection on "Forecast Model A - Data" (it previously held the
# tab-selected view with no explicit selection; it now has an explicit
# selection but is no longer the selected tab).
$wsAData.Range("J2:J101").Select()

# Make "Forecast Model B - Data" the active sheet/tab, preserving its
# existing J2:J101 selection.
$wsBData.Activate()
$wsBData.Range("J2:J101").Select()
